$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header row (copy format from H1 so both reuse the same style).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data values for columns I and J, rows 2-17.
$data = @(
    @(2, 6, 6),
    @(3, 7, 7),
    @(4, 8, 8),
    @(5, 5, 5),
    @(6, 6, 7),
    @(7, 6, 6),
    @(8, 8, 8),
    @(9, 5, 6),
    @(10, 7, 8),
    @(11, 8, 8),
    @(12, 4, 5),
    @(13, 8, 9),
    @(14, 9, 9),
    @(15, 9, 9),
    @(16, 7, 7),
    @(17, 6, 6)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
